# Generate Report for Handback
# Updates the localization-status report after a successful handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - The "Latest Handback DateTime" is refreshed to the handback run time
#  - The stale "Error Detail" (handback-not-latest warning) is cleared now
#    that the handback is in sync
#  - A couple of columns are widened/narrowed to fit the new text

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("K2").Value = "2016-09-02 01:00:42"
$zhcn.Range("P2").Value = $null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Range("K2").Value = "2016-09-02 01:00:49"
$dede.Range("P2").Value = $null

# ---------------------------------------------------------------------------
# Overview sheet mirrors each language's Status in its own column
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText

# ---------------------------------------------------------------------------
# Column width touch-ups (Status/handback columns got wider text, the old
# long Error Detail column can shrink back down)
# ---------------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.16666667
$overview.Columns.Item(6).ColumnWidth = 29.16666667

$zhcn.Columns.Item(3).ColumnWidth = 29.16666667
$zhcn.Columns.Item(16).ColumnWidth = 12.83333333

$dede.Columns.Item(3).ColumnWidth = 29.16666667
$dede.Columns.Item(16).ColumnWidth = 12.83333333
